$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Weekly Sales": append a new weekly row (row 19) continuing the
# 7-day cadence of column A, with the corresponding y value in column B.
# ---------------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Sales")
$wsWeekly.Range("A19").Value = 45662.99999999999
$wsWeekly.Range("A19").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsWeekly.Range("B19").Value = 76

# ---------------------------------------------------------------------------
# Sheet "Merged (Optional)": shuffle several Daily_PO_Qty (column C) values
# around within their date groups (Penalty Reward System experiment), and
# append a new row 63 mirroring the new weekly row above.
# ---------------------------------------------------------------------------
$wsMerged = $wb.Worksheets.Item("Merged (Optional)")

$wsMerged.Range("C30").Value = 12
$wsMerged.Range("C31").Value = 4
$wsMerged.Range("C34").Value = 4
$wsMerged.Range("C35").Value = 8
$wsMerged.Range("C37").Value = 4
$wsMerged.Range("C38").Value = 12
$wsMerged.Range("C40").Value = 8
$wsMerged.Range("C44").Value = 8
$wsMerged.Range("C45").Value = 4
$wsMerged.Range("C48").Value = 4
$wsMerged.Range("C49").Value = 76
$wsMerged.Range("C50").Value = 4
$wsMerged.Range("C51").Value = 12
$wsMerged.Range("C52").Value = 76
$wsMerged.Range("C53").Value = 12

$wsMerged.Range("A63").Value = 45662.99999999999
$wsMerged.Range("A63").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsMerged.Range("B63").Value = 76
$wsMerged.Range("C63").Value = 0
